# inverter_config.xlsx — "update load file excel set on-off"
#
# 1) Stations sheet: the old "B2 / INV-10" row (row 11) was removed from the
#    station list, so every row below it shifts up by one (dimension goes
#    from A1:F134 to A1:F133).
# 2) Control_Scenarios sheet: the OFF/ON scenario counts (column C) were
#    updated to reflect the new inverter groupings.

$wb = $excel.ActiveWorkbook

# --- 1) Stations: delete the obsolete B2/INV-10 row, shifting rows up ---
$ws1 = $wb.Worksheets.Item("Stations")
$ws1.Rows("11:11").Delete()

# --- 2) Control_Scenarios: refresh the scenario inverter counts ---
$ws2 = $wb.Worksheets.Item("Control_Scenarios")

# Values are stored as text (matches the sheet's existing numberStoredAsText
# cells), so format the range as text before writing the new counts.
$ws2.Range("C2:C8").NumberFormat = "@"

$ws2.Range("C2").Value = "9"   # B3R1 OFF
$ws2.Range("C3").Value = "10"  # B4R2 OFF
$ws2.Range("C4").Value = "10"  # B5R2 OFF
$ws2.Range("C5").Value = "5"   # B8   OFF
$ws2.Range("C6").Value = "9"   # B3R1 ON
$ws2.Range("C7").Value = "10"  # B4R2 ON
$ws2.Range("C8").Value = "10"  # B5R2 ON
# C9 (B8 ON) stays "5" - unchanged

# Drop back to the workbook's default (unstyled) cell style so only the
# value/type changes, not the formatting, same as the rest of the sheet.
$ws2.Range("C2:C8").Style = "Normal"
